# Apply roster updates to rows 3-16 (columns A:C) on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3,  "Luka Doncic",         "PG,SG",    "Dallas Mavericks"),
    @(4,  "DeMar DeRozan",       "SF,PF",    "Sacramento Kings"),
    @(5,  "Brook Lopez",         "C",        "Milwaukee Bucks"),
    @(6,  "Jusuf Nurkic",        "C",        "Phoenix Suns"),
    @(7,  "Mikal Bridges",       "SG,SF,PF", "New York Knicks"),
    @(8,  "Evan Mobley",         "PF,C",     "Cleveland Cavaliers"),
    @(10, "Amen Thompson",       "SG,SF",    "Houston Rockets"),
    @(11, "Michael Porter Jr.",  "SF,PF",    "Denver Nuggets"),
    @(12, "Bilal Coulibaly",     "SG,SF",    "Washington Wizards"),
    @(13, "Bennedict Mathurin",  "SG,SF",    "Indiana Pacers"),
    @(14, "Tyler Herro",         "PG,SG",    "Miami Heat"),
    @(15, "Josh Giddey",         "PG,SG,SF", "Chicago Bulls"),
    @(16, "Nikola Vucevic",      "PF,C",     "Chicago Bulls")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
